$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1) ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"
$ws.Range("F1:K1").Style = $ws.Range("B1").Style

# --- Data rows (rows 2-7): company, name, owner, property_category ---
$ws.Range("B2").Value = "國泰人壽"
$ws.Range("C2").Value = "得意還本終身"
$ws.Range("D2").Value = "陳明文"
$ws.Range("E2").Value = "insurance"

$ws.Range("B3").Value = "國泰人壽"
$ws.Range("C3").Value = "新富貴保本投資鏈結型保險第7期"
$ws.Range("D3").Value = "廖素惠"
$ws.Range("E3").Value = "insurance"

$ws.Range("B4").Value = "國泰人壽"
$ws.Range("C4").Value = "創世變額萬能壽險"
$ws.Range("D4").Value = "廖素惠"
$ws.Range("E4").Value = "insurance"

$ws.Range("B5").Value = "國泰人壽"
$ws.Range("C5").Value = "創世變額萬能壽險"
$ws.Range("D5").Value = "廖素惠"
$ws.Range("E5").Value = "insurance"

$ws.Range("B6").Value = "富邦人壽"
$ws.Range("C6").Value = "安泰還本終身壽險"
$ws.Range("D6").Value = "廖素惠"
$ws.Range("E6").Value = "insurance"

$ws.Range("B7").Value = "中國信託人壽"
$ws.Range("C7").Value = "年年沛616美元還本終身壽險"
$ws.Range("D7").Value = "陳明文"
$ws.Range("E7").Value = "insurance"

# --- Common trailer columns F-K for every data row ---
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 6).Value = "normal"
    $ws.Cells.Item($r, 7).Value = "2013-11-20"
    $ws.Cells.Item($r, 8).Value = "陳明文"
    $ws.Cells.Item($r, 9).Value = 828
    $ws.Cells.Item($r, 10).Value = "tmp581f1"
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($r, 1).Value
    $ws.Range($ws.Cells.Item($r, 6), $ws.Cells.Item($r, 11)).Style = $ws.Range("B2").Style
}
